# Refresh the cryptocurrency Price/Volume(1h) snapshot columns (D, E) to
# the latest scrape, including the Monero / InternetComputer(DFINITY) row
# swap at rows 26-27. All Price values are stored as text (several look
# like plain numbers, e.g. '328.03', so a leading apostrophe forces text
# the way typing them into Excel would).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.991.29"
$ws.Range("E2").Value = "  +1.64%  "

$ws.Range("D3").Value = "'1.945.25"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.41%  "

$ws.Range("D5").Value = "'327.67"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").Value = "'0.4847"
$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'0.4095"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").Value = "'0.08222"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").Value = "'1.017"
$ws.Range("E10").Value = "  -0.92%  "

$ws.Range("D11").Value = "'24.02"
$ws.Range("E11").Value = "  +1.04%  "

$ws.Range("D12").Value = "'1.965.78"
$ws.Range("E12").Value = "  +3.65%  "

$ws.Range("D13").Value = "'6.105"
$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("D14").Value = "'7.336"
$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("D15").Value = "'91.85"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").Value = "'0.06867"
$ws.Range("E16").Value = "  +1.21%  "

$ws.Range("D17").Value = "'1.010"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").Value = "'0.00001039"
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").Value = "'17.83"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").Value = "'1.008"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").Value = "'29.977.21"
$ws.Range("E21").Value = "  +1.46%  "

$ws.Range("D22").Value = "'5.694"
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("D23").Value = "'11.95"
$ws.Range("E23").Value = "  +1.53%  "

$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("D25").Value = "'2.188.77"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'6.573"
$ws.Range("E26").Value = "  -2.28%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'156.77"
$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("D28").Value = "'20.10"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").Value = "'2.117"
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("D30").Value = "'121.23"
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("D32").Value = "'0.09627"
$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("D33").Value = "'5.627"
$ws.Range("E33").Value = "  +1.70%  "

$ws.Range("D34").Value = "'1.422"
$ws.Range("E34").Value = "  +1.99%  "

$ws.Range("D35").Value = "'3.553"
$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("D36").Value = "'0.06545"
$ws.Range("E36").Value = "  +6.33%  "

$ws.Range("D37").Value = "'0.02300"
$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("D38").Value = "'1.214"
$ws.Range("E38").Value = "  +2.72%  "

$ws.Range("D39").Value = "'0.5961"
$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("D40").Value = "'10.75"
$ws.Range("E40").Value = "  -0.78%  "

$ws.Range("D41").Value = "'7.961"
$ws.Range("E41").Value = "  -1.25%  "

$ws.Range("D42").Value = "'2.543"
$ws.Range("E42").Value = "  +4.86%  "

$ws.Range("D43").Value = "'0.1854"
$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("D44").Value = "'12.50"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "'1.246"
$ws.Range("E45").Value = "  -2.61%  "

$ws.Range("D46").Value = "'0.07559"
$ws.Range("E46").Value = "  -0.60%  "

$ws.Range("D47").Value = "'0.5582"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("D48").Value = "'1.988"
$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D49").Value = "'118.06"
$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("D50").Value = "'2.434"
$ws.Range("E50").Value = "  -0.22%  "

$ws.Range("D51").Value = "'72.58"
$ws.Range("E51").Value = "  -0.43%  "

